$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracked price-history row (2026-02-07 scrape).
# The source column values are stored as text in this workbook (shared
# strings), not as native numbers/dates, so force a text number format
# before writing the values to avoid Excel auto-converting "2026-02-07"
# into a date serial or "48800000"/"0" into native numbers.
$row = 36
$rngAddress = "A" + $row + ":D" + $row
$rng = $ws.Range($rngAddress)
$rng.NumberFormat = "@"

$ws.Range("A" + $row).Value = "2026-02-07"
$ws.Range("B" + $row).Value = "48800000"
$ws.Range("C" + $row).Value = "0"
$ws.Range("D" + $row).Value = "0"

# Restore the cells to the workbook's default (unstyled) look, matching
# every other row in the sheet.
$rng.Style = "Normal"
